{"js": "// Add two new \"place\" entries (label paragraph + hyperlink paragraph) right\n// after the existing \"Karwendelparkplatz Rum:\" hyperlink paragraph:\n//   Gasthof Koreth:\n//     https://www.innsbruck.info/fr/hebergements/hebergements/unterkunft/gasthof-koreth-innsbruck.html\n//   Schie\u00dfstand Arzl:\n//     https://www.almenrausch.at/uploads/tx_wctrip/DSC_9005_15959_01.jpg\n\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph whose hyperlink run text is the Karwendelparkplatz\n// Rum hyperlink URL -- the new content is inserted directly after it.\nconst anchorText =\n  \"https://www.innsbruck.info/infrastruktur/detail/infrastruktur/parkplatz-alpenpark-karwendel-rum.html\";\nlet anchorIndex = -1;\nfor (let i = 0; i < paras.items.length; i++) {\n  if (paras.items[i].text === anchorText) {\n    anchorIndex = i;\n    break;\n  }\n}\nif (anchorIndex === -1) {\n  throw new Error(\"Could not find the Karwendelparkplatz Rum hyperlink paragraph\");\n}\n\n// Helper: build a minimal single-document OOXML wrapper for insertOoxml().\nfunction wrapOoxml(bodyXml) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:r=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships\">' +\n    \"<w:body>\" +\n    bodyXml +\n    \"</w:body>\" +\n    \"</w:document>\" +\n    \"</pkg:xmlData></pkg:part></pkg:package>\"\n  );\n}\n\n// Inserts a plain label paragraph (e.g. \"Gasthof Koreth:\") with the middle\n// word wrapped in proofErr spell-check markers, right after `afterPara`.\n// Returns the inserted paragraph's text so the caller can re-find it.\nasync function insertLabelParagraph(afterPara, prefix, spellWord, suffix) {\n  const bodyXml =\n    \"<w:p>\" +\n    '<w:r><w:t xml:space=\"preserve\">' +\n    prefix +\n    \"</w:t></w:r>\" +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    \"<w:r><w:t>\" +\n    spellWord +\n    \"</w:t></w:r>\" +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    \"<w:r><w:t>\" +\n    suffix +\n    \"</w:t></w:r>\" +\n    \"</w:p>\";\n  const placeholder = afterPara.insertParagraph(\"\", Word.InsertLocation.after);\n  placeholder.insertOoxml(wrapOoxml(bodyXml), Word.InsertLocation.replace);\n  await context.sync();\n  return prefix + spellWord + suffix;\n}\n\n// Inserts a paragraph containing a single hyperlink run (styled with the\n// built-in \"Hyperlink\" character style), right after `afterPara`. Using the\n// Range.hyperlink setter (rather than raw OOXML) so Word mints a proper\n// external relationship + applies the Hyperlink style.\nasync function insertHyperlinkParagraph(afterPara, url) {\n  const p = afterPara.insertParagraph(url, Word.InsertLocation.after);\n  const rng = p.getRange();\n  rng.hyperlink = url;\n  await context.sync();\n  return url;\n}\n\n// Re-locate a paragraph by its exact text (paragraphs collection indices\n// shift after each insertion, so we look the anchor back up each time).\nasync function findParagraphByText(text) {\n  paras.load(\"items/text\");\n  await context.sync();\n  for (let i = 0; i < paras.items.length; i++) {\n    if (paras.items[i].text === text) {\n      return paras.items[i];\n    }\n  }\n  throw new Error(\"Could not re-locate paragraph with text: \" + text);\n}\n\nlet cursor = paras.items[anchorIndex];\n\nconst gasthofText = await insertLabelParagraph(cursor, \"Gasthof \", \"Koreth\", \":\");\ncursor = await findParagraphByText(gasthofText);\n\nawait insertHyperlinkParagraph(\n  cursor,\n  \"https://www.innsbruck.info/fr/hebergements/hebergements/unterkunft/gasthof-koreth-innsbruck.html\"\n);\ncursor = await findParagraphByText(\n  \"https://www.innsbruck.info/fr/hebergements/hebergements/unterkunft/gasthof-koreth-innsbruck.html\"\n);\n\nconst schiessstandText = await insertLabelParagraph(cursor, \"Schie\u00dfstand \", \"Arzl\", \":\");\ncursor = await findParagraphByText(schiessstandText);\n\nawait insertHyperlinkParagraph(\n  cursor,\n  \"https://www.almenrausch.at/uploads/tx_wctrip/DSC_9005_15959_01.jpg\"\n);\n\nawait context.sync();\n", "ps1": "# Add two new \"place\" entries (label paragraph + hyperlink paragraph) right\n# after the existing \"Karwendelparkplatz Rum:\" hyperlink paragraph:\n#   Gasthof Koreth:\n#     https://www.innsbruck.info/fr/hebergements/hebergements/unterkunft/gasthof-koreth-innsbruck.html\n#   Schie\u00dfstand Arzl:\n#     https://www.almenrausch.at/uploads/tx_wctrip/DSC_9005_15959_01.jpg\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph whose text is the Karwendelparkplatz Rum hyperlink URL\n# -- the new content is inserted directly after it.\n$anchorText = \"https://www.innsbruck.info/infrastruktur/detail/infrastruktur/parkplatz-alpenpark-karwendel-rum.html\"\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    # Range.Text includes the trailing paragraph mark (\"`r\"), so trim it\n    # before comparing against the plain anchor text.\n    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd(\"`r\") -eq $anchorText) {\n        $anchorIndex = $i\n        break\n    }\n}\nif ($anchorIndex -eq -1) {\n    throw \"Could not find the Karwendelparkplatz Rum hyperlink paragraph\"\n}\n\n# Inserts a plain label paragraph (e.g. \"Gasthof Koreth:\") with the middle\n# word wrapped in proofErr spell-check markers, right after $afterPara.\n# Returns the newly inserted Paragraph object.\nfunction Insert-LabelParagraph($afterPara, $prefix, $spellWord, $suffix) {\n    $insertAt = $afterPara.Index\n    [void]$afterPara.Range.InsertParagraphAfter()\n    $newPara = $d.Paragraphs.Item($insertAt + 1)\n    $xml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:r><w:t xml:space=\"preserve\">' + $prefix + '</w:t></w:r>' +\n        '<w:proofErr w:type=\"spellStart\"/>' +\n        '<w:r><w:t>' + $spellWord + '</w:t></w:r>' +\n        '<w:proofErr w:type=\"spellEnd\"/>' +\n        '<w:r><w:t>' + $suffix + '</w:t></w:r>' +\n        '</w:p>'\n    $newPara.Range.InsertXML($xml)\n    return $d.Paragraphs.Item($insertAt + 1)\n}\n\n# Inserts a paragraph containing a single hyperlink run (styled with the\n# built-in \"Hyperlink\" character style), right after $afterPara. Uses\n# Hyperlinks.Add (rather than raw XML) so Word mints a proper external\n# relationship + applies the Hyperlink style; the range is shrunk by one\n# character first so the trailing paragraph mark isn't swept into the new\n# hyperlink run.\nfunction Insert-HyperlinkParagraph($afterPara, $url) {\n    $insertAt = $afterPara.Index\n    [void]$afterPara.Range.InsertParagraphAfter()\n    $newPara = $d.Paragraphs.Item($insertAt + 1)\n    $newPara.Range.Text = $url\n    $newPara2 = $d.Paragraphs.Item($insertAt + 1)\n    $hrng = $newPara2.Range\n    $hrng.MoveEnd(1, -1)\n    [void]$d.Hyperlinks.Add($hrng, $url)\n    return $d.Paragraphs.Item($insertAt + 1)\n}\n\n$anchor = $d.Paragraphs.Item($anchorIndex)\n\n$p1 = Insert-LabelParagraph $anchor \"Gasthof \" \"Koreth\" \":\"\n$p2 = Insert-HyperlinkParagraph $p1 \"https://www.innsbruck.info/fr/hebergements/hebergements/unterkunft/gasthof-koreth-innsbruck.html\"\n$p3 = Insert-LabelParagraph $p2 \"Schie\u00dfstand \" \"Arzl\" \":\"\n$p4 = Insert-HyperlinkParagraph $p3 \"https://www.almenrausch.at/uploads/tx_wctrip/DSC_9005_15959_01.jpg\"\n"}
